$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------------
# Old layout (columns A..S):
#   A Norm, Typ | B Varumärke | C Artikelbenämning | D GVM | E Artikelnummer
#   F Typbeteckning | G Ritningsnummer | H Position | I Beteckning
#   J Kompletterande Information övrigt | K Ref annan | L Historiskt Varumärke
#   M Historiskt inköpsreferens | N Enhet | O Förpackning | P SSG-notering
#   Q (empty) | R E-nummer | S RSK-nummer
#
# New layout (columns A..I) keeps only:
#   Varumärke, Artikelbenämning, GVM, Artikelnummer, Typbeteckning, Enhet,
#   SSG-notering, E-nummer, RSK-nummer
# ----------------------------------------------------------------------------------

# Remove the unused/empty column between "SSG-notering" (P) and "E-nummer" (R)
$ws.Range("Q1").EntireColumn.Delete() | Out-Null

# Remove "Förpackning"
$ws.Range("O1").EntireColumn.Delete() | Out-Null

# Remove "Ritningsnummer", "Position", "Beteckning", "Kompletterande Information övrigt",
# "Ref annan", "Historiskt Varumärke", "Historiskt inköpsreferens"
$ws.Range("G1:M1").EntireColumn.Delete() | Out-Null

# Remove "Norm, Typ"
$ws.Range("A1").EntireColumn.Delete() | Out-Null

# Remove the autofilter (and its embedded sort state) that used to span the old columns
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Update the hidden _FilterDatabase defined name to reflect the new column extent
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$I`$21217"
    }
}

# Widen the "SSG-notering" column (now column G) to fit its content
$ws.Columns.Item(7).ColumnWidth = 28.6667

# Update the view: selection now targets column H ("E-nummer") instead of the old Q,
# and the previously frozen/scrolled-to top-left cell (J1) is no longer relevant
$ws.Range("H1:H1048576").Select() | Out-Null

$wb.Save()
